# Updated Code & Feature File
# Add a new "ServiceCleanup" worksheet at the end of the workbook, populate
# it with the EVC/OVC and UNI service cleanup headers/fields, and switch the
# active/selected tab to EVC_Disconnect.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ServiceCleanup"

# --- Populate the cell content ---------------------------------------------
$newSheet.Range("A1").Value = "EVC/OVC End Point Cleanup"
$newSheet.Range("A2").Value = "IdentifierId"
$newSheet.Range("A3").Value = "CorrelationId"
$newSheet.Range("A5").Value = "UNI Service Cleanup"
$newSheet.Range("A6").Value = "UNI Service Id"
$newSheet.Range("A7").Value = "CustomerName"
$newSheet.Range("A8").Value = "CustomerNo"

# --- Bold the two section headers ------------------------------------------
$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A5").Font.Bold = $true

# --- Match column widths used by the rest of the workbook ------------------
$newSheet.Columns.Item(1).ColumnWidth = 27.1666666666667
$newSheet.Columns.Item(2).ColumnWidth = 21.3072916666667

# --- Leave the selection on the last populated cell -------------------------
[void]$newSheet.Range("A8").Select()

# --- Restore EVC_Disconnect as the active/selected sheet -------------------
$evcDisc = $wb.Worksheets.Item("EVC_Disconnect")
[void]$evcDisc.Activate()
